$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated coin-ranking snapshot (symbol list refresh).
# The source sheet stores every value (prices, % deltas, names, links)
# as literal text, so we force NumberFormat "@" (Text) on each cell
# before writing it -- otherwise Excel would auto-convert values like
# "37.27" or "-0.04%" into real numbers/percentages.
$updates = [ordered]@{
    "D2" = "307.87"
    "D3" = "37.27"
    "E3" = "-0.04%"
    "D4" = "5.120"
    "E4" = "-0.25%"
    "D5" = "0.07829"
    "E5" = "0.76%"
    "B6" = "GateToken"
    "C6" = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
    "D6" = "4.394"
    "E6" = "-0.13%"
    "B7" = "KuCoinToken"
    "C7" = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
    "D7" = "8.252"
    "E7" = "0.73%"
    "B8" = "FTXToken"
    "C8" = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
    "D8" = "1.881"
    "E8" = "0.35%"
    "B9" = "BTSEToken"
    "C9" = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
    "D9" = "2.943"
    "E9" = "2.28%"
    "B10" = "MXToken"
    "C10" = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
    "D10" = "0.9209"
    "E10" = "0.30%"
    "B11" = "LiechtensteinCryptoassetsExchange"
    "C11" = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
    "D11" = "0.1081"
    "E11" = "-9.14%"
    "B12" = "WazirX"
    "C12" = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
    "D12" = "0.1888"
    "E12" = "-0.35%"
    "B13" = "MandalaExchangeToken"
    "C13" = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
    "D13" = "0.08848"
    "E13" = "-6.54%"
    "B14" = "BitrueCoin"
    "C14" = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
    "D14" = "0.03311"
    "E14" = "-2.85%"
    "B15" = "BitMartToken"
    "C15" = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
    "D15" = "0.09583"
    "E15" = "-1.12%"
    "B16" = "BitForexToken"
    "C16" = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
    "D16" = "0.001380"
    "E16" = "0.51%"
    "B17" = "TigerCash"
    "C17" = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
    "D17" = "0.005821"
    "E17" = "0.05%"
    "B18" = "LEO"
    "C18" = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
    "D18" = "3.402"
    "E18" = "-3.80%"
    "E19" = "0.71%"
    "D20" = "6.303"
    "E20" = "19.90%"
    "D21" = "0.1303"
    "E21" = "2.89%"
    "E22" = "-6.73%"
    "D23" = "0.04366"
    "E23" = "0.95%"
    "E24" = "-0.30%"
    "E25" = "0.30%"
    "D26" = "0.0001400"
    "E26" = "7.55%"
    "E27" = "-98.10%"
    "D39" = "0.02174"
    "E39" = "4.83%"
    "E40" = "0.50%"
    "D41" = "0.007550"
    "E41" = "-1.71%"
    "D42" = "0.1353"
    "E42" = "0.53%"
    "D43" = "0.008651"
    "E43" = "-11.92%"
    "D44" = "0.002068"
    "E44" = "-4.81%"
    "D45" = "0.007901"
    "E45" = "-9.51%"
    "D46" = "0.00006518"
    "E46" = "-2.89%"
    "E47" = "-0.12%"
    "D48" = "0.003293"
    "E48" = "12.15%"
    "E49" = "-16.58%"
    "E50" = "-0.12%"
    "E51" = "-0.12%"
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
}
